$d = $word.ActiveDocument

# 1. "54645" -> "6463"  (NÚMERO ÚNICO DE EMPLEADO value)
$d.Content.Find.Execute("54645", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6463", 2) | Out-Null

# 2. "30 DE NOVIEMBRE DE 2025" -> "26 DE JUNIO DE 2026"
$d.Content.Find.Execute("30 DE NOVIEMBRE DE 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "26 DE JUNIO DE 2026", 2) | Out-Null
